$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value, taken from the diff.
# D column holds prices, E column holds 1h volume/change percentages;
# both are stored as plain text (inline strings) in the sheet, so we
# force a text NumberFormat before writing to stop Excel from coercing
# values such as "1.00" or "316.89" into numbers, then clear the
# temporary formatting so the cell style is left untouched.
$updates = [ordered]@{
    "D2" = "47.830.23"
    "E2" = "  -1.01%  "
    "D3" = "2.481.92"
    "E3" = "  -1.68%  "
    "E4" = "  -0.03%  "
    "D5" = "316.89"
    "E5" = "  -1.69%  "
    "D6" = "104.06"
    "E6" = "  -5.28%  "
    "D7" = "0.517"
    "E7" = "  -2.81%  "
    "D8" = "1.00"
    "E8" = "  -0.02%  "
    "D9" = "0.534"
    "E9" = "  -3.08%  "
    "D10" = "38.73"
    "E10" = "  -4.50%  "
    "D11" = "20.34"
    "E11" = "  -0.44%  "
    "E12" = "  -3.02%  "
    "E13" = "  +0.49%  "
    "D14" = "7.01"
    "E14" = "  -3.60%  "
    "D15" = "2.872.03"
    "E15" = "  -1.71%  "
    "D16" = "2.372.91"
    "E16" = "  -6.00%  "
    "D17" = "0.822"
    "E17" = "  -3.77%  "
    "D18" = "47.792.95"
    "E18" = "  -0.75%  "
    "D19" = "12.65"
    "E19" = "  -5.10%  "
    "D20" = "2.89"
    "E20" = "  +7.60%  "
    "D21" = "6.50"
    "E21" = "  -1.80%  "
    "E22" = "  -2.54%  "
    "D23" = "278.69"
    "E23" = "  +3.33%  "
    "D24" = "70.50"
    "E24" = "  -2.02%  "
    "D25" = "2.49"
    "E25" = "  -3.64%  "
    "E26" = "  -0.20%  "
    "D27" = "25.56"
    "E27" = "  -1.80%  "
    "D28" = "2.21"
    "E28" = "  -1.83%  "
    "D29" = "9.58"
    "E29" = "  -5.57%  "
    "D30" = "0.137"
    "E30" = "  -4.96%  "
    "D31" = "34.33"
    "E31" = "  -3.59%  "
    "D32" = "48.98"
    "E32" = "  -1.35%  "
    "E33" = "  -0.08%  "
    "D34" = "18.91"
    "E34" = "  -4.18%  "
    "D35" = "5.23"
    "E35" = "  -2.91%  "
    "D36" = "0.0766"
    "E36" = "  -2.48%  "
    "E37" = "  -3.00%  "
    "E38" = "  -4.86%  "
    "E39" = "  -5.51%  "
    "D40" = "121.61"
    "E40" = "  -0.28%  "
    "E41" = "  -1.51%  "
    "D42" = "2.22"
    "E42" = "  +0.47%  "
    "D43" = "21.30"
    "E43" = "  -2.58%  "
    "D44" = "0.0297"
    "E44" = "  -1.63%  "
    "D45" = "1.981.98"
    "E45" = "  -2.17%  "
    "D46" = "3.10"
    "E46" = "  -1.76%  "
    "D47" = "1.90"
    "E47" = "  +0.42%  "
    "D48" = "2.04"
    "E48" = "  +0.08%  "
    "D49" = "8.88"
    "E49" = "  -2.63%  "
    "D50" = "5.06"
    "E50" = "  -2.93%  "
    "D51" = "78.70"
    "E51" = "  -0.96%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
